{"js": "// Office.js (Word JavaScript API) script.\n// Applies the English -> Swahili subtitle translation edits described\n// by the diff. Each entry is an exact (before, after) text pair; we\n// locate the text with Body.search (exact, case-sensitive, single\n// occurrence) and swap it in place with Range.insertText(..., \"Replace\")\n// so existing run formatting (bold/italic/fonts/etc.) on the run is\n// preserved.\n\nconst pairs = [\n  [\"Prisoners and candies - subtitles:\", \"Wafungwa na peremende - manukuu:\"],\n  [\"**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino\", \"**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino\"],\n  [\"[Music]\", \"[Muziki]\"],\n  [\"four bright mathematicians are taken into\", \"wanahisabati wanne mkali wanachukuliwa\"],\n  [\"custody and put in jail because they tried\", \"chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu\"],\n  [\"to convince an old lady that the Goedel's\", \"kumshawishi mwanamke mzee kuwa Goedel's\"],\n  [\"incompleteness theorems are true. Every\", \"nadharia za kutokamilika ni kweli. Kila\"],\n  [\"mathematician has his own cell that we\", \"mtaalamu wa hisabati ana kiini chake ambacho sisi\"],\n  [\"can enumerate with a number from 1 to 4.\", \"inaweza kuhesabu na nambari kutoka 1 hadi 4.\"],\n  [\"before entering the cell a certain\", \"kabla ya kuingia kwenye seli fulani\"],\n  [\"number of candies greater than \", \"idadi ya peremende kubwa kuliko \"],\n  [\"e: OR \", \"e: AU \"],\n  [\"EQUAL TO\", \"SAWA NA\"],\n  [\" 1 is\", \" 1 ni\"],\n  [\"given to every mathematician and they\", \"wanapewa kila mtaalamu wa hisabati na wao\"],\n  [\"are told they have 11 candies in total.\", \"wanaambiwa wana peremende 11 kwa jumla.\"],\n  [\"but everyone knows only his number of\", \"lakini kila mtu anajua idadi yake tu\"],\n  [\"candies and the total. 1 and is not\", \"pipi na jumla. 1 na sio\"],\n  [\"allowed to ask for the others number.\", \"kuruhusiwa kuuliza nambari zingine.\"],\n  [\"then the first mathematician asks the\", \"kisha mwanahisabati wa kwanza anauliza\"],\n  [\"second: 'number 2 do you know if you\", \"pili: 'namba 2 unajua kama wewe\"],\n  [\"have more candies than me?' the second\", \"kuwa na peremende nyingi kuliko mimi?' ya pili\"],\n  [\"mathematician answers he doesn't. Then he\", \"mwanahisabati anajibu hana. Kisha yeye\"],\n  [\"asks to number 3: 'do you know if you have\", \"anauliza kwa nambari 3: 'unajua kama unayo\"],\n  [\"more candy than me?'\", \"pipi zaidi kuliko mimi?'\"],\n  [\"the third mathematician answers: 'no I'm\", \"mwanahisabati wa tatu anajibu: 'hapana niko\"],\n  [\"sorry I don't'. At this point the fourth\", \"samahani sifanyi'. Katika hatua hii ya nne\"],\n  [\"mathematician says: 'hey guys you know\", \"mtaalamu wa hisabati anasema: 'jamani mnafahamu\"],\n  [\"what, I know exactly how many candies\", \"nini, najua hasa pipi ngapi\"],\n  [\"everyone has here'. Surprisingly even the\", \"kila mtu ana hapa'. Cha kushangaza hata\"],\n  [\"other three mathematicians say that now\", \"wanahisabati wengine watatu wanasema hivyo sasa\"],\n  [\"they know how many candies everyone has\", \"wanajua kila mtu ana pipi ngapi\"],\n  [\"so the question is: can you figure out\", \"kwa hivyo swali ni: unaweza kujua\"],\n  [\"the number of candies every prisoner has\", \"idadi ya pipi kila mfungwa ana\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + before);\n  }\n\n  // Replace every occurrence (there are exactly two identical\n  // \"[Music]\" cues in this subtitle file; every other string is\n  // unique in the document).\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// The diff also inserts an extra, standalone space between the\n// closing \")\" of the \"(errata corrige: ...)\" aside and the following\n// \" 1 ni\" text, so the rendered line goes from \"...) 1 is\" to\n// \"...)  1 ni\" (note the double space). The preceding loop already\n// turned \" 1 is\" into \" 1 ni\" as its own run (preserving that run's\n// non-italic formatting), so just insert a single space immediately\n// before it, which extends the (italic) \")\" run instead of touching\n// the (non-italic) \" 1 ni\" run.\nconst beforeNi = context.document.body.search(\" 1 ni\", { matchCase: true });\nbeforeNi.load(\"items\");\nawait context.sync();\nif (beforeNi.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for ' 1 ni'\");\n}\nbeforeNi.items[0].insertText(\" \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the English -> Swahili subtitle translation edits described\n# by the diff. Each entry is an exact (before, after) text pair; we\n# use Find/Replace (wdReplaceAll) against the whole document Range so\n# existing run formatting (bold/italic/fonts/etc.) on the matched text\n# is preserved, same as a manual Ctrl+H in Word.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"Prisoners and candies - subtitles:\", \"Wafungwa na peremende - manukuu:\"),\n  @(\"**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino\", \"**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino\"),\n  @(\"[Music]\", \"[Muziki]\"),\n  @(\"four bright mathematicians are taken into\", \"wanahisabati wanne mkali wanachukuliwa\"),\n  @(\"custody and put in jail because they tried\", \"chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu\"),\n  @(\"to convince an old lady that the Goedel's\", \"kumshawishi mwanamke mzee kuwa Goedel's\"),\n  @(\"incompleteness theorems are true. Every\", \"nadharia za kutokamilika ni kweli. Kila\"),\n  @(\"mathematician has his own cell that we\", \"mtaalamu wa hisabati ana kiini chake ambacho sisi\"),\n  @(\"can enumerate with a number from 1 to 4.\", \"inaweza kuhesabu na nambari kutoka 1 hadi 4.\"),\n  @(\"before entering the cell a certain\", \"kabla ya kuingia kwenye seli fulani\"),\n  @(\"number of candies greater than \", \"idadi ya peremende kubwa kuliko \"),\n  @(\"e: OR \", \"e: AU \"),\n  @(\"EQUAL TO\", \"SAWA NA\"),\n  @(\" 1 is\", \" 1 ni\"),\n  @(\"given to every mathematician and they\", \"wanapewa kila mtaalamu wa hisabati na wao\"),\n  @(\"are told they have 11 candies in total.\", \"wanaambiwa wana peremende 11 kwa jumla.\"),\n  @(\"but everyone knows only his number of\", \"lakini kila mtu anajua idadi yake tu\"),\n  @(\"candies and the total. 1 and is not\", \"pipi na jumla. 1 na sio\"),\n  @(\"allowed to ask for the others number.\", \"kuruhusiwa kuuliza nambari zingine.\"),\n  @(\"then the first mathematician asks the\", \"kisha mwanahisabati wa kwanza anauliza\"),\n  @(\"second: 'number 2 do you know if you\", \"pili: 'namba 2 unajua kama wewe\"),\n  @(\"have more candies than me?' the second\", \"kuwa na peremende nyingi kuliko mimi?' ya pili\"),\n  @(\"mathematician answers he doesn't. Then he\", \"mwanahisabati anajibu hana. Kisha yeye\"),\n  @(\"asks to number 3: 'do you know if you have\", \"anauliza kwa nambari 3: 'unajua kama unayo\"),\n  @(\"more candy than me?'\", \"pipi zaidi kuliko mimi?'\"),\n  @(\"the third mathematician answers: 'no I'm\", \"mwanahisabati wa tatu anajibu: 'hapana niko\"),\n  @(\"sorry I don't'. At this point the fourth\", \"samahani sifanyi'. Katika hatua hii ya nne\"),\n  @(\"mathematician says: 'hey guys you know\", \"mtaalamu wa hisabati anasema: 'jamani mnafahamu\"),\n  @(\"what, I know exactly how many candies\", \"nini, najua hasa pipi ngapi\"),\n  @(\"everyone has here'. Surprisingly even the\", \"kila mtu ana hapa'. Cha kushangaza hata\"),\n  @(\"other three mathematicians say that now\", \"wanahisabati wengine watatu wanasema hivyo sasa\"),\n  @(\"they know how many candies everyone has\", \"wanajua kila mtu ana pipi ngapi\"),\n  @(\"so the question is: can you figure out\", \"kwa hivyo swali ni: unaweza kujua\"),\n  @(\"the number of candies every prisoner has\", \"idadi ya pipi kila mfungwa ana\"),\n)\n\nforeach ($pair in $pairs) {\n  $searchText = $pair[0]\n  $replaceText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $searchText\n  $find.Replacement.Text = $replaceText\n\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n  #         ReplaceWith, Replace)\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $result = $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n  if (-not $result) {\n    throw \"Text not found: $searchText\"\n  }\n}\n\n# The diff also inserts an extra, standalone space between the closing\n# \")\" of the \"(errata corrige: ...)\" aside and the following \" 1 ni\"\n# text, so the rendered line goes from \"...) 1 is\" to \"...)  1 ni\"\n# (note the double space). The loop above already turned \" 1 is\" into\n# \" 1 ni\" as its own run (preserving that run's non-italic\n# formatting), so just insert a single space immediately before it,\n# which extends the (italic) \")\" run instead of touching the\n# (non-italic) \" 1 ni\" run.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \" 1 ni\"\n$found2 = $find2.Execute()\nif (-not $found2) {\n  throw \"Text not found: ' 1 ni'\"\n}\n$niRange = $find2.Parent\n$insertPoint = $d.Range($niRange.Start, $niRange.Start)\n$insertPoint.InsertBefore(\" \")\n"}
